# FormatKM.xlsx edit script
# Adds two custom number formats (K/M abbreviations) applied via formulas in
# columns B & D (referencing column A), plus text columns C & E holding the
# expected formatted text of B & D for comparison/testing purposes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Custom number formats -------------------------------------------------
$fmtK = '[>999999]#,,"M";[>999]#,"K";#'
$fmtM = '[>999999]#.000,,"M";[>999]#.000,"K";#.000'

# --- Column layout -----------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 13.140625

# --- Text values that each row should display as in columns B (K fmt) and D (M fmt)
$kText = @{
    2  = "1"
    3  = "10"
    4  = "102"
    5  = "102"
    6  = "1K"
    7  = "10K"
    8  = "102K"
    9  = "1M"
    10 = "10M"
    11 = "102M"
    12 = "1021M"
}
$mText = @{
    2  = "1.020"
    3  = "10.200"
    4  = "102.000"
    5  = "102.102"
    6  = "1.021K"
    7  = "10.210K"
    8  = "102.102K"
    9  = "1.021M"
    10 = "10.210M"
    11 = "102.102M"
    12 = "1021.021M"
}

for ($r = 2; $r -le 12; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Formula = "=A$r"
    $bCell.NumberFormat = $fmtK

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Formula = "=A$r"
    $dCell.NumberFormat = $fmtM
}

for ($r = 2; $r -le 12; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value = "'" + $kText[$r]

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = "'" + $mText[$r]
}

# --- Page setup (A4, portrait) ---------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---------------------------------------------------------
$ws.Range("F8").Select()

$wb.Save()
